$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.311.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.930.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.83%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7207'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3297'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.88'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06925'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8024'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08066'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.929.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.410'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.300.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008356'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.802'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.184.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.864'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.719'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.402'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.42%  '
$ws.Range('E29').Value = '  -12.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.554'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.51%  '
$ws.Range('E31').Value = '  -1.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.403'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.189'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.29%  '
$ws.Range('E34').Value = '  -2.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.220'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7397'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.741'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.29%  '
$ws.Range('E38').Value = '  -2.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.829'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.588'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '78.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4461'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.990'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.45%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8360'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.794'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.294'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4079'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.93%  '
